$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14), shifting the existing
# "Late" / "heading" / "Outstanding" columns one place to the right
# (N->O, O->P, P->Q). Use the width of the preceding column (M, "In Advance")
# for the newly inserted column, matching what Excel does on a manual
# column insert.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, with cell R8 selected.
$ws.Activate()
$ws.Range("R8").Select()
